$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.192.33"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.201.99"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.200.07"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").Value = "3.727.76"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "66.350.96"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("D18").Value = "3.203.81"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "507.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +39.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  -4.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "501.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "0.0₃0774"
$ws.Range("E39").Value = "  +15.44%  "
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0421"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "2.898.86"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
